$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Kode Fixed Income (generated bond code) from OBL00107 to OBL00108
$ws.Range("M2").Value = "OBL00108"

# Update the big preparation text block to reference the new generated bond code
$text = "Username : 31246;`nPassword : bni1234;`nRole : 20/21 - Analis Investasi/Asisten Investasi;`nKode Fixed Income : OBL00108;`nNama Fixed Income : Obligasi III Mitra Jaya Tahun 2019"
$ws.Range("F2").Value = $text

# Row 2 is a touch shorter now that the text fits differently
$ws.Range("A2:X2").RowHeight = 105

# Update the saved view/selection state - scroll so column B is the
# left-most visible column, then select G2
$win = $wb.Windows.Item(1)
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("G2").Select() | Out-Null
